# "Add Data Sink Option"
# Replace the old Email/First Name/Last Name/Avatar columns with a new
# Id/Name/Job table populated from a data sink (reqres.in "users" / "jobs"
# style data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two columns (D, E) that are no longer part of the layout.
$ws.Range("D1:E1").Clear()

# Rewrite the remaining headers.
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Job"

# The raw "job" payload is an (empty) XML element coming back from the
# data source for the first three rows.
$jobTag = '<job xmlns="https://reqres.in/api/users"/>'

$data = @(
  @("470", "Engineer", $jobTag),
  @("398", "Teacher",  $jobTag),
  @("8",   "Doctor",   $jobTag),
  @("318", "chathuri", "Engineer"),
  @("542", "buddhika", "Teacher"),
  @("865", "gunapala", "Doctor")
)

# The Id column holds numeric-looking values that must stay text, so force
# the range to Text formatting before writing the values, then restore the
# default style once the text values are in place.
$idRange = $ws.Range("A2:A7")
$idRange.NumberFormat = "@"

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r++
}

$idRange.Style = "Normal"

# Match the author's final selection in the sheet.
$ws.Range("G7").Select()
